$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.174.46'
$ws.Range('E2').Value = '  -5.80%  '
$ws.Range('D3').Value = '2.453.78'
$ws.Range('E3').Value = '  -8.50%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.70%  '
$ws.Range('D9').Value = '2.471.54'
$ws.Range('E9').Value = '  -8.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0994'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.76%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('D14').Value = '2.897.06'
$ws.Range('E14').Value = '  -8.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.13%  '
$ws.Range('D16').Value = '59.067.71'
$ws.Range('E16').Value = '  -5.85%  '
$ws.Range('E17').Value = '  -5.45%  '
$ws.Range('D18').Value = '2.503.00'
$ws.Range('E18').Value = '  -6.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.969'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.450'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -11.79%  '
$ws.Range('E26').Value = '  -5.09%  '
$ws.Range('E27').Value = '  -2.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.79%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.47%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.30%  '
$ws.Range('E31').Value = '  -5.57%  '
$ws.Range('D32').Value = '0.0₃0772'
$ws.Range('E32').Value = '  -8.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.997'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '157.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.41'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.45'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.72'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '314.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.31'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.835'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.65%  '
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.594'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.78%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0528'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0934'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0230'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.60%  '
